$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data as described in the commit diff.
# Each cell is forced to plain text formatting before/after the write so that
# numeric-looking strings (e.g. "0.911", "0.0000283") are preserved verbatim as
# text rather than being auto-converted to numbers by Excel, while avoiding any
# residual "text" number-format style being left behind on the cell.

function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.311.44"
Set-TextValue $ws.Range("E2") "  +1.42%  "
Set-TextValue $ws.Range("D3") "3.407.68"
Set-TextValue $ws.Range("E3") "  +1.49%  "
Set-TextValue $ws.Range("E4") "  +0.09%  "
Set-TextValue $ws.Range("D5") "582.44"
Set-TextValue $ws.Range("E5") "  -0.46%  "
Set-TextValue $ws.Range("D6") "179.07"
Set-TextValue $ws.Range("E6") "  +1.05%  "
Set-TextValue $ws.Range("E7") "  +0.00%  "
Set-TextValue $ws.Range("E8") "  +0.46%  "
Set-TextValue $ws.Range("E9") "  +8.02%  "
Set-TextValue $ws.Range("E10") "  +0.78%  "
Set-TextValue $ws.Range("D11") "48.46"
Set-TextValue $ws.Range("E11") "  +0.97%  "
Set-TextValue $ws.Range("D12") "0.0000283"
Set-TextValue $ws.Range("E12") "  +3.13%  "
Set-TextValue $ws.Range("D13") "682.42"
Set-TextValue $ws.Range("E13") "  -0.60%  "
Set-TextValue $ws.Range("B14") "Polkadot"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "8.65"
Set-TextValue $ws.Range("E14") "  +2.31%  "
Set-TextValue $ws.Range("B15") "WrappedliquidstakedEther2.0"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D15") "3.957.93"
Set-TextValue $ws.Range("E15") "  +1.28%  "
Set-TextValue $ws.Range("D16") "69.498.54"
Set-TextValue $ws.Range("E16") "  +1.62%  "
Set-TextValue $ws.Range("D17") "3.404.33"
Set-TextValue $ws.Range("E17") "  +1.57%  "
Set-TextValue $ws.Range("E18") "  +0.55%  "
Set-TextValue $ws.Range("D19") "17.71"
Set-TextValue $ws.Range("E19") "  +1.16%  "
Set-TextValue $ws.Range("D20") "11.32"
Set-TextValue $ws.Range("E20") "  +0.90%  "
Set-TextValue $ws.Range("D21") "0.911"
Set-TextValue $ws.Range("E21") "  +1.70%  "
Set-TextValue $ws.Range("E22") "  -1.92%  "
Set-TextValue $ws.Range("D23") "17.08"
Set-TextValue $ws.Range("E23") "  +0.72%  "
Set-TextValue $ws.Range("D24") "100.78"
Set-TextValue $ws.Range("E24") "  +0.59%  "
Set-TextValue $ws.Range("E25") "  -0.44%  "
Set-TextValue $ws.Range("E26") "  -0.05%  "
Set-TextValue $ws.Range("D27") "9.72"
Set-TextValue $ws.Range("E27") "  +2.06%  "
Set-TextValue $ws.Range("D28") "33.55"
Set-TextValue $ws.Range("E28") "  +1.58%  "
Set-TextValue $ws.Range("D29") "8.76"
Set-TextValue $ws.Range("E29") "  +2.65%  "
Set-TextValue $ws.Range("E30") "  -1.36%  "
Set-TextValue $ws.Range("D31") "3.77"
Set-TextValue $ws.Range("E31") "  +12.20%  "
Set-TextValue $ws.Range("D32") "559.09"
Set-TextValue $ws.Range("E32") "  +1.43%  "
Set-TextValue $ws.Range("E33") "  -0.64%  "
Set-TextValue $ws.Range("E34") "  -0.01%  "
Set-TextValue $ws.Range("D35") "57.99"
Set-TextValue $ws.Range("E35") "  -0.23%  "
Set-TextValue $ws.Range("E36") "  +0.11%  "
Set-TextValue $ws.Range("D37") "3.611.50"
Set-TextValue $ws.Range("E37") "  -2.86%  "
Set-TextValue $ws.Range("D38") "0.140"
Set-TextValue $ws.Range("E38") "  +0.56%  "
Set-TextValue $ws.Range("D39") "35.24"
Set-TextValue $ws.Range("E39") "  +1.48%  "
Set-TextValue $ws.Range("D40") "0.0₃0745"
Set-TextValue $ws.Range("E40") "  +10.55%  "
Set-TextValue $ws.Range("D41") "3.29"
Set-TextValue $ws.Range("E41") "  +3.54%  "
Set-TextValue $ws.Range("E42") "  +3.17%  "
Set-TextValue $ws.Range("E43") "  +3.26%  "
Set-TextValue $ws.Range("E44") "  +3.36%  "
Set-TextValue $ws.Range("E45") "  +0.13%  "
Set-TextValue $ws.Range("E46") "  +1.16%  "
Set-TextValue $ws.Range("E47") "  +0.24%  "
Set-TextValue $ws.Range("D48") "1.41"
Set-TextValue $ws.Range("E48") "  +4.65%  "
Set-TextValue $ws.Range("E49") "  -0.20%  "
Set-TextValue $ws.Range("D50") "131.49"
Set-TextValue $ws.Range("E50") "  -0.14%  "
Set-TextValue $ws.Range("E51") "  +3.45%  "
